# Append a second row of data to the worksheet (row 2), mirroring the
# header row 1 layout: سؤال 8 | خيار 1..5 | 5 | شرح | اسم المادة | اسم الفصل | اسم التصنيف
#
# Numeric-looking values ("1","2","3","4","5") must be stored as TEXT
# (not numbers), so they are entered with a leading apostrophe, which is
# the standard Excel convention for forcing text/"number stored as text".
# Purely alphabetic values are entered normally (they are text already).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "سؤال"
$ws.Range("B2").Value = "'1"
$ws.Range("C2").Value = "'2"
$ws.Range("D2").Value = "'3"
$ws.Range("E2").Value = "'4"
$ws.Range("F2").Value = "'5"
$ws.Range("G2").Value = "'5"
$ws.Range("H2").Value = "واو"
$ws.Range("I2").Value = "علوم"
$ws.Range("J2").Value = "الاول"
$ws.Range("K2").Value = "جلد"
